$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table used to carry an extra leading "group_id" index column (A) that
# duplicated the real group_id values already living in column B. Drop that
# column so B..E (group_id, treatment, color, independant_variables) slide
# left into A..D - this also naturally carries row 2's blank trailing cell
# (old E2) along into the new D2, and keeps the header's bold/boxed style
# lined up over the shifted columns.
$ws.Columns("A").Delete()

# Add the new trailing "label" column (E), populated with the short
# treatment name per group.
$ws.Range("E1").Value = "label"
$ws.Range("E2").Value = "vehicles"
$ws.Range("E3").Value = "MDL"
$ws.Range("E4").Value = "TCB2"
$ws.Range("E5").Value = "TCB2+MDL"

# Match the header style (bold, bordered, centered) already used by A1:D1
# on the new E1 header cell.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
